$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 4
$ws.Range("F5").Value = 377
$ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202410/SPyugcNX1729569622422.jpeg"
$ws.Range("F7").Value = 3910
$ws.Range("F9").Value = 773
$ws.Range("F10").Value = 2323
$ws.Range("F16").Value = 2249
$ws.Range("F17").Value = 323
$ws.Range("F18").Value = 27
$ws.Range("F20").Value = 346
$ws.Range("F22").Value = 45
$ws.Range("F23").Value = 277

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 50
$ws.Range("F5").Value = 37
$ws.Range("F7").Value = 133
$ws.Range("F10").Value = 99
$ws.Range("F12").Value = 6
$ws.Range("F22").Value = 64

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6400
$ws.Range("F5").Value = 341

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 6400
$ws.Range("F5").Value = 341
$ws.Range("F6").Value = 50
$ws.Range("F7").Value = 50
$ws.Range("F8").Value = 4
$ws.Range("F12").Value = 377
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202410/SPyugcNX1729569622422.jpeg"
$ws.Range("F14").Value = 37
$ws.Range("F18").Value = 3910
$ws.Range("F19").Value = 133
$ws.Range("F23").Value = 99
$ws.Range("F24").Value = 773
$ws.Range("F25").Value = 2323
$ws.Range("F32").Value = 6
$ws.Range("F34").Value = 2249
$ws.Range("F35").Value = 323
$ws.Range("F38").Value = 27
$ws.Range("F40").Value = 346
$ws.Range("F42").Value = 45
$ws.Range("F49").Value = 64
$ws.Range("F50").Value = 277
